# lead-bulk-template.xlsx: add a "propertyType" column between "budget" and
# "flatType" (new column E), pushing flatType/areaKey/remark one column to
# the right (E->F, F->G, G->H), and add its list-validation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at E; this shifts existing data, column widths and
# data validations for columns E,F,G one place to the right automatically.
$ws.Columns.Item(5).Insert()

# Header text for the newly inserted column.
$ws.Range("E1").Value = "propertyType"

# Restore the width of the now-split column block: D keeps width 15, the
# new E should be 20, F keeps 15 (same as the old E/flatType column).
$ws.Columns.Item(5).ColumnWidth = 19.166666666666668

# Re-create the validations for the shifted columns (F = flatType,
# G = areaKey) so they end up listed in column order, then add the new
# validation for E (propertyType) in between D and F.
$ws.Range("F2:F5000").Validation.Delete()
$ws.Range("G2:G5000").Validation.Delete()

$ws.Range("E2:E5000").Validation.Add(3, 1, 1, '"Standalone house,Apartment,Gated community,Independent house,Villa,PG / Co-living,Plot / Land"')
$ws.Range("E2:E5000").Validation.ShowInput = $false
$ws.Range("E2:E5000").Validation.ShowError = $false

$ws.Range("F2:F5000").Validation.Add(3, 1, 1, '"1RK,1BHK,2BHK,3BHK,4BHK,Villa,Penthouse"')
$ws.Range("F2:F5000").Validation.ShowInput = $false
$ws.Range("F2:F5000").Validation.ShowError = $false

$ws.Range("G2:G5000").Validation.Add(3, 1, 1, '"Whitefield,Indiranagar,Koramangala,Bengaluru,Jayanagar,Banashankari,Basaveshwaranagar,Bheemanahalli,Bommanahalli,Chikkalasandra,Dasarahalli,Domlur,Electronic City,Frazer Town,Girinagar,Gokula,Gopalapuram,Hanumanthanagar,HBR Layout,Hebbal,Hoysala,HSR Layout,Ittamadu,JP Nagar,Jyothinagar,Kammanahalli,Kaval Byrasandra,Kodichikkanahalli,Kommadi,Kundalahalli,Lingrajapuram,Mahadevapura,Malleswaram,Marathahalli,Mathikere,Mico Layout,Mookambika,Nagavara,Nagawara,Nagarathpet,Nandini Layout,Nayandahalli,Old Airport Road,Peenya,Prithviraj Road,RMV Extension,Sadashivnagar,Sahakarnagar,Sanjaynagar,Sarjapur Road,Seshadripuram,Shantinagar,Shivaji Nagar,Soladevanahalli,Subramanyanagar"')
$ws.Range("G2:G5000").Validation.ShowInput = $false
$ws.Range("G2:G5000").Validation.ShowError = $false

Write-Host "lead-bulk-template: propertyType column inserted"
